# Insert 3 new data rows right before the current row 32 (old rows 32-65
# shift down to 35-68). This matches the diff: dimension grows from
# A1:R65 to A1:R68, and all previously-existing rows reappear unmodified
# three rows further down.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(32).Resize(3).Insert()

# Constant values shared by every data row in this sheet.
$mercado   = "Agrícola del Norte S.A. de Arica"
$region    = "Arica y Parinacota"
$codreg    = 15
$catId     = 100112008
$categoria = "Coliflor"
$variedad  = "Sin especificar"
$unidad    = "`$/unidad"
$origen    = "Región de Arica y Parinacota"
$clasif    = "Hortaliza"

# New rows 32-34: three quality grades reported for fecha 2021-09-28 (44467).
$newRows = @(
    @{ Row = 32; Fecha = 44467; Calidad = "Primera"; Vol = 600;  PMin = 1000; PMax = 1200; PProm = 1100; PKg = 1100 },
    @{ Row = 33; Fecha = 44467; Calidad = "Segunda"; Vol = 800;  PMin = 700;  PMax = 800;  PProm = 750;  PKg = 750  },
    @{ Row = 34; Fecha = 44467; Calidad = "Tercera"; Vol = 800;  PMin = 400;  PMax = 500;  PProm = 450;  PKg = 450  }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = 1
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $r.Fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $catId
    $ws.Cells.Item($row, 7).Value  = $categoria
    $ws.Cells.Item($row, 8).Value  = $variedad
    $ws.Cells.Item($row, 9).Value  = $r.Calidad
    $ws.Cells.Item($row, 10).Value = $r.Vol
    $ws.Cells.Item($row, 11).Value = $r.PMin
    $ws.Cells.Item($row, 12).Value = $r.PMax
    $ws.Cells.Item($row, 13).Value = $r.PProm
    $ws.Cells.Item($row, 14).Value = $unidad
    $ws.Cells.Item($row, 15).Value = $origen
    $ws.Cells.Item($row, 16).Value = $r.PKg
    $ws.Cells.Item($row, 17).Value = 1
    $ws.Cells.Item($row, 18).Value = $clasif
}
